$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A159").Value = 29
